$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

$ws.Range('D2').Value = '28.630.41'
$ws.Range('E2').Value = '  +0.77%  '
$ws.Range('D3').Value = '1.800.83'
$ws.Range('E3').Value = '  -0.89%  '
Set-TextValue $ws.Range('D5') '315.92'
$ws.Range('E5').Value = '  -0.63%  '
Set-TextValue $ws.Range('D6') '1.002'
$ws.Range('E6').Value = '  +0.34%  '
Set-TextValue $ws.Range('D7') '0.5309'
$ws.Range('E7').Value = '  -7.73%  '
Set-TextValue $ws.Range('D8') '0.3767'
$ws.Range('E8').Value = '  -2.53%  '
Set-TextValue $ws.Range('D9') '42.48'
$ws.Range('E9').Value = '  -1.97%  '
Set-TextValue $ws.Range('D10') '0.07487'
$ws.Range('E10').Value = '  -1.94%  '
$ws.Range('E11').Value = '  -2.49%  '
$ws.Range('E12').Value = '  +0.41%  '
Set-TextValue $ws.Range('D13') '20.66'
$ws.Range('E13').Value = '  -2.98%  '
Set-TextValue $ws.Range('D14') '6.142'
$ws.Range('E14').Value = '  -1.98%  '
Set-TextValue $ws.Range('D15') '7.337'
$ws.Range('E15').Value = '  +0.39%  '
$ws.Range('D16').Value = '1.794.77'
$ws.Range('E16').Value = '  -1.16%  '
Set-TextValue $ws.Range('D17') '90.18'
$ws.Range('E17').Value = '  -2.45%  '
Set-TextValue $ws.Range('D18') '0.00001063'
$ws.Range('E18').Value = '  -1.73%  '
Set-TextValue $ws.Range('D19') '0.06463'
$ws.Range('E19').Value = '  -0.83%  '
$ws.Range('E20').Value = '  +0.28%  '
$ws.Range('E21').Value = '  -0.74%  '
Set-TextValue $ws.Range('D22') '5.890'
$ws.Range('E22').Value = '  -1.80%  '
$ws.Range('D23').Value = '28.634.23'
$ws.Range('E23').Value = '  +0.72%  '
Set-TextValue $ws.Range('D24') '11.07'
$ws.Range('E24').Value = '  -2.74%  '
Set-TextValue $ws.Range('D25') '2.090'
$ws.Range('E25').Value = '  -0.41%  '
Set-TextValue $ws.Range('D26') '159.74'
$ws.Range('E26').Value = '  +1.39%  '
$ws.Range('E27').Value = '  -2.30%  '
$ws.Range('D28').Value = '1.998.86'
$ws.Range('E28').Value = '  -1.37%  '
Set-TextValue $ws.Range('D29') '2.340'
$ws.Range('E29').Value = '  -3.30%  '
Set-TextValue $ws.Range('D30') '122.38'
$ws.Range('E30').Value = '  -1.10%  '
Set-TextValue $ws.Range('D31') '1.100'
$ws.Range('E31').Value = '  -5.72%  '
Set-TextValue $ws.Range('D32') '0.1050'
$ws.Range('E32').Value = '  -1.08%  '
Set-TextValue $ws.Range('D33') '3.697'
$ws.Range('E33').Value = '  +2.05%  '
Set-TextValue $ws.Range('D34') '5.630'
$ws.Range('E34').Value = '  -2.76%  '
Set-TextValue $ws.Range('D35') '0.2246'
$ws.Range('E35').Value = '  +3.78%  '
Set-TextValue $ws.Range('D36') '0.06404'
$ws.Range('E36').Value = '  +5.15%  '
Set-TextValue $ws.Range('D37') '0.02305'
$ws.Range('E37').Value = '  -0.58%  '
Set-TextValue $ws.Range('D38') '8.802'
$ws.Range('E38').Value = '  -0.82%  '
Set-TextValue $ws.Range('D39') '5.035'
$ws.Range('E39').Value = '  -0.33%  '
Set-TextValue $ws.Range('D40') '1.212'
$ws.Range('E40').Value = '  +4.66%  '
Set-TextValue $ws.Range('D41') '11.23'
$ws.Range('E41').Value = '  -4.27%  '
Set-TextValue $ws.Range('D42') '0.6192'
$ws.Range('E42').Value = '  -3.75%  '
$ws.Range('E43').Value = '  +0.28%  '
Set-TextValue $ws.Range('D44') '1.409'
$ws.Range('E44').Value = '  +2.14%  '
Set-TextValue $ws.Range('D45') '13.33'
$ws.Range('E45').Value = '  -1.31%  '
$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range('D46') '3.687'
$ws.Range('E46').Value = '  -0.49%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue $ws.Range('D47') '0.5840'
$ws.Range('E47').Value = '  -2.52%  '
Set-TextValue $ws.Range('D48') '125.60'
$ws.Range('E48').Value = '  +2.73%  '
Set-TextValue $ws.Range('D49') '1.937'
$ws.Range('E49').Value = '  -0.41%  '
$ws.Range('E50').Value = '  -0.04%  '
Set-TextValue $ws.Range('D51') '0.06886'
$ws.Range('E51').Value = '  +0.48%  '
